$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.275.77'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '2.613.25'
$ws.Range("E3").Value = '  +0.50%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''584.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.83%  '
$ws.Range("D6").Value = '''143.24'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.71%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").Value = '''0.597'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = '''6.49'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("E10").Value = '  -0.56%  '
$ws.Range("E11").Value = '  +2.19%  '
$ws.Range("D12").Value = '''0.155'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.92%  '
$ws.Range("D13").Value = '3.073.66'
$ws.Range("E13").Value = '  -0.41%  '
$ws.Range("D14").Value = '''25.02'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.67%  '
$ws.Range("D15").Value = '60.275.18'
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("D16").Value = '''0.0000140'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("D17").Value = '2.616.19'
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").Value = '''11.46'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.66%  '
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("D20").Value = '''347.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.38%  '
$ws.Range("E21").Value = '  -2.24%  '
$ws.Range("E22").Value = '  -0.52%  '
$ws.Range("D23").Value = '''0.531'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.37%  '
$ws.Range("D24").Value = '''63.68'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.92%  '
$ws.Range("D25").Value = '''0.998'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("E26").Value = '  +0.31%  '
$ws.Range("D27").Value = '''8.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.30%  '
$ws.Range("D28").Value = '''1.94'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.09%  '
$ws.Range("E29").Value = '  +1.24%  '
$ws.Range("E30").Value = '  +2.43%  '
$ws.Range("D31").Value = '''168.60'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.06%  '
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("E33").Value = '  +0.16%  '
$ws.Range("E34").Value = '  +6.20%  '
$ws.Range("E35").Value = '  +8.66%  '
$ws.Range("D36").Value = '''4.28'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.10%  '
$ws.Range("E37").Value = '  +2.29%  '
$ws.Range("D38").Value = '''319.92'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.95%  '
$ws.Range("D39").Value = '''38.45'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.93%  '
$ws.Range("E40").Value = '  +4.19%  '
$ws.Range("D41").Value = '''0.851'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.31%  '
$ws.Range("D42").Value = '''135.52'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.42%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = '''20.04'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.92%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").Value = '''0.0991'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.60%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '''0.999'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("E46").Value = '  +0.85%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '''5.04'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.48%  '
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").Value = '''0.0553'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.35%  '
$ws.Range("D49").Value = '''20.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.90%  '
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").Value = '''10.76'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.56%  '
